# Journal update: mark Day 23 (Nmap) task Done and Day 24 (TryHackMe Vulnversity) as In Progress.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "Hands-on tools: Nmap, Wireshark, basic scans" -> Status: Done, In Progress? unchecked, Done? checked
$ws.Range("C10").Value = "Done"
$ws.Range("D10").Value = "☐"
$ws.Range("E10").Value = "☑"

# Row 11: "Do beginner labs on TryHackMe / HackTheBox" -> Status: In Progress, In Progress? checked
$ws.Range("C11").Value = "In Progress"
$ws.Range("D11").Value = "☑"

# Update print orientation to portrait for the sheet
$ws.PageSetup.Orientation = 1

# Move the active selection to B11, matching the author's last edit position
$ws.Range("B11").Select()

$wb.Save()
